$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Custom Statuses")

$rows = @(
    @('Using Eyespore as an aggro card', 'Durga'),
    @('Ign: 40anb swr', 'Durga'),
    @('Sunshroom is destroying the meta, nerf Miner', 'FrozenFighter'),
    @('With Reincarnation, VRH and Fig, we could be able to get a Reliable Molekale..... stares at camera', 'StingRay'),
    @('When your Molekale spawns a Pineclone', 'Frozen Fighter'),
    @('turn 1 black hole turn 2 double cheese', 'Huyn'),
    @('rose is so OP why won''t my 4x5 HG conjure pile work on her???', 'StingRay'),
    @('Fighting for Corndog', 'Waltuh'),
    @('Which hero can play OTK Cryo Yeti? HG? No. BF? No. GS? Yes.', 'connorthepeach'),
    @('Helping Vengeful Cyborg look for his lost toupee', 'Pasta')
)

$startRow = 299
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
